$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Sheet 1
$ws1.Range("F4").Value = 841
$ws1.Range("F7").Value = 9613
$ws1.Range("G7").Value = "已售罄"
$ws1.Range("F10").Value = 712
$ws1.Range("F11").Value = 2098
$ws1.Range("F12").Value = 49
$ws1.Range("F13").Value = 1623
$ws1.Range("F14").Value = 2724
$ws1.Range("F15").Value = 135
$ws1.Range("F16").Value = 4061
$ws1.Range("F17").Value = 333
$ws1.Range("F18").Value = 160
$ws1.Range("F20").Value = 218
$ws1.Range("F22").Value = 31
$ws1.Range("F23").Value = 80
$ws1.Range("F24").Value = 79
$ws1.Range("F26").Value = 3824
$ws1.Range("F28").Value = 3321
$ws1.Range("F30").Value = 199
$ws1.Range("F32").Value = 4332
$ws1.Range("F34").Value = 276
$ws1.Range("F35").Value = 393
$ws1.Range("F36").Value = 257

# Sheet 2
$ws2.Range("F3").Value = 22

# Sheet 3
$ws3.Range("F3").Value = 995

# Sheet 4
$ws4.Range("F4").Value = 995
$ws4.Range("F6").Value = 841
$ws4.Range("F9").Value = 9613
$ws4.Range("G9").Value = "已售罄"
$ws4.Range("F12").Value = 712
$ws4.Range("F13").Value = 2098
$ws4.Range("F14").Value = 49
$ws4.Range("F15").Value = 1623
$ws4.Range("F17").Value = 2724
$ws4.Range("F18").Value = 135
$ws4.Range("F19").Value = 4061
$ws4.Range("F20").Value = 333
$ws4.Range("F21").Value = 160
$ws4.Range("F23").Value = 218
$ws4.Range("F25").Value = 31
$ws4.Range("F26").Value = 22
$ws4.Range("F27").Value = 80
$ws4.Range("F28").Value = 80
$ws4.Range("F30").Value = 3825
$ws4.Range("F32").Value = 3321
$ws4.Range("F34").Value = 199
$ws4.Range("F35").Value = 490
$ws4.Range("F36").Value = 4332
$ws4.Range("F38").Value = 276
$ws4.Range("F39").Value = 393
$ws4.Range("F40").Value = 257
